$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) C2, C3, C4: "67%" -> "73%" (kept as literal text, same style)
# ------------------------------------------------------------------
foreach ($r in 2..4) {
    $target = "C$r"
    $src = "B$r"
    $ws.Range($target).NumberFormat = "@"
    $ws.Range($target).Value = "73%"
    # restore the original cell formatting (style index) that the
    # NumberFormat tweak perturbed, by re-applying the format from the
    # neighbouring cell that already carries the right style.
    $ws.Range($src).Copy()
    $ws.Range($target).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Grow the log with new sensor readings.
#    Original rows 14 & 15 (no explicit row height, style 11 = "last
#    row" look) become regular data rows (style 10, ht=15), two more
#    regular rows (16, 17) are appended, and a new final row (18)
#    takes over the "last row" styling that 14/15 used to have.
# ------------------------------------------------------------------

# 2a) Capture the "last row" look (from row 15, still untouched) onto
#     the brand-new row 18 before we restyle 14/15.
$ws.Range("A15:C15").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)

# 2b) Turn rows 14 and 15 into normal data rows, matching rows 10-13.
foreach ($r in 14, 15) {
    $ws.Range("A10:C10").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
    $ws.Rows($r).RowHeight = 15
}

# 2c) Add two new normal data rows (16, 17).
foreach ($r in 16, 17) {
    $ws.Range("A10:C10").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
    $ws.Rows($r).RowHeight = 15
}
$excel.CutCopyMode = 0

# 2d) Fill in the values for rows 14-18.
$ws.Range("A14").Value = "26/04/2025 14:00:15"
$ws.Range("B14").Value = "25°"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "67%"
$ws.Range("B14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("A15").Value = "26/04/2025 14:00:31"
$ws.Range("B15").Value = "25°"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "67%"
$ws.Range("B15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("A16").Value = "26/04/2025 14:24:18"
$ws.Range("B16").Value = "25°"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "67%"
$ws.Range("B16").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("A17").Value = "26/04/2025 14:25:16"
$ws.Range("B17").Value = "25°"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "67%"
$ws.Range("B17").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("A18").Value = "26/04/2025 14:30:08"
$ws.Range("B18").Value = "25°"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "65%"
$ws.Range("B18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Move the active selection to H14, as in the edited workbook.
# ------------------------------------------------------------------
[void]$ws.Range("H14").Select()
